$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# This workbook has 3 sheets: "Overview", "zh-cn", "de-de".
# Each sheet has a row for each tracked file, identified by column A
# (a hyperlinked file name / guid). We need to insert two brand new rows
# - one for 27f14fa0-3a00-4acd-b6db-b7bf1cb1f21b.md
# - one for 91b753c9-5d1c-49cf-8bd3-2fbeee9ef532.md
# right before the existing "a6289566-8408-4fec-9a8c-04180367b81e.md" row,
# on every sheet. The simplest reliable way (since Hyperlinks do not shift
# when rows are inserted in this host) is to:
#   1) duplicate the "a6289566" row twice (copy/insert) so the three rows
#      keep identical shape/styles,
#   2) overwrite the text in the two new rows with the correct values,
#   3) rebuild the full Hyperlinks collection from scratch in row order.
# -------------------------------------------------------------------------

function Rebuild-Overview-Sheet($ws) {
    # Row 3 ("a6289566...") is duplicated twice, pushing it down to row 5,
    # and ".localization-config" down to row 6.
    $ws.Rows.Item(3).Copy()
    $ws.Rows.Item(3).Insert()
    $ws.Rows.Item(3).Copy()
    $ws.Rows.Item(3).Insert()

    # Fill in the new row 3 (27f14fa0)
    $ws.Range("A3").Value = "27f14fa0-3a00-4acd-b6db-b7bf1cb1f21b.md"
    $ws.Range("B3").Value = "Ready for handoff"
    $ws.Range("C3").Value = "Ready for handoff"

    # Fill in the new row 4 (91b753c9)
    $ws.Range("A4").Value = "91b753c9-5d1c-49cf-8bd3-2fbeee9ef532.md"
    $ws.Range("B4").Value = "Ready for handoff"
    $ws.Range("C4").Value = "Ready for handoff"

    # Rebuild hyperlinks top to bottom so rIds line up with row order.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1b714050bd19de299f2bc9b207fc0bdde29e62a9/e2e/9ad15719-6ce1-48cc-8569-036f14eacdc0.md", [Type]::Missing, [Type]::Missing, "9ad15719-6ce1-48cc-8569-036f14eacdc0.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/29f8223291f41775c48fc5f81140839eeff1258d/e2e/27f14fa0-3a00-4acd-b6db-b7bf1cb1f21b.md", [Type]::Missing, [Type]::Missing, "27f14fa0-3a00-4acd-b6db-b7bf1cb1f21b.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/29f8223291f41775c48fc5f81140839eeff1258d/e2e/91b753c9-5d1c-49cf-8bd3-2fbeee9ef532.md", [Type]::Missing, [Type]::Missing, "91b753c9-5d1c-49cf-8bd3-2fbeee9ef532.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/29f8223291f41775c48fc5f81140839eeff1258d/e2e/a6289566-8408-4fec-9a8c-04180367b81e.md", [Type]::Missing, [Type]::Missing, "a6289566-8408-4fec-9a8c-04180367b81e.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/1b714050bd19de299f2bc9b207fc0bdde29e62a9/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null
}

function Rebuild-Locale-Sheet($ws, $locale, $c2Date1, $c2Date2, $c3Date, $c4Date, $c5Date) {
    # Row 3 ("a6289566...") is duplicated twice, pushing it down to row 5,
    # and ".localization-config" down to row 6.
    $ws.Rows.Item(3).Copy()
    $ws.Rows.Item(3).Insert()
    $ws.Rows.Item(3).Copy()
    $ws.Rows.Item(3).Insert()

    # Fill in the new row 3 (27f14fa0)
    $ws.Range("A3").Value = "27f14fa0-3a00-4acd-b6db-b7bf1cb1f21b.md"
    $ws.Range("B3").Value = "Ready for handoff"
    $ws.Range("C3").Value = "27f14fa0-3a00-4acd-b6db-b7bf1cb1f21b.2b242c02220495576e019c268293a2629917dac4.$locale.xlf"
    $ws.Range("D3").Value = $c3Date
    $ws.Range("G3").Value = "0001-01-01 00:00:00"
    $ws.Range("H3").Value = "Include"

    # Fill in the new row 4 (91b753c9)
    $ws.Range("A4").Value = "91b753c9-5d1c-49cf-8bd3-2fbeee9ef532.md"
    $ws.Range("B4").Value = "Ready for handoff"
    $ws.Range("C4").Value = "91b753c9-5d1c-49cf-8bd3-2fbeee9ef532.5dd33dce7d1a10de4f7d8751457d7bbec4f4e1ab.$locale.xlf"
    $ws.Range("D4").Value = $c4Date
    $ws.Range("G4").Value = "0001-01-01 00:00:00"
    $ws.Range("H4").Value = "Include"

    # Rebuild hyperlinks top to bottom so rIds line up with row order.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1b714050bd19de299f2bc9b207fc0bdde29e62a9/e2e/9ad15719-6ce1-48cc-8569-036f14eacdc0.md", [Type]::Missing, [Type]::Missing, "9ad15719-6ce1-48cc-8569-036f14eacdc0.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5fe2827f2b618c4dbb3dd4f6b978994cc35b86d4/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/ht/9ad15719-6ce1-48cc-8569-036f14eacdc0.79be59d77f61d35d36ff06509bb7aaec296dacce.$locale.xlf", [Type]::Missing, [Type]::Missing, "9ad15719-6ce1-48cc-8569-036f14eacdc0.79be59d77f61d35d36ff06509bb7aaec296dacce.$locale.xlf") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.$locale/blob/87e613b29c933314248b3e34ece6faf9c018e5ac/e2e/9ad15719-6ce1-48cc-8569-036f14eacdc0.md", [Type]::Missing, [Type]::Missing, "9ad15719-6ce1-48cc-8569-036f14eacdc0.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e9b98f292261f6a73c0f455b7df23ae8b55bd0bb/ol-handback/OpenLocalizationTestOrg/oltest.$locale/ci/ht/9ad15719-6ce1-48cc-8569-036f14eacdc0.79be59d77f61d35d36ff06509bb7aaec296dacce.$locale.xlf", [Type]::Missing, [Type]::Missing, "9ad15719-6ce1-48cc-8569-036f14eacdc0.79be59d77f61d35d36ff06509bb7aaec296dacce.$locale.xlf") | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/29f8223291f41775c48fc5f81140839eeff1258d/e2e/27f14fa0-3a00-4acd-b6db-b7bf1cb1f21b.md", [Type]::Missing, [Type]::Missing, "27f14fa0-3a00-4acd-b6db-b7bf1cb1f21b.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/16eb241e7716f10c06e8c8193ec79918d1a355bd/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/ht/27f14fa0-3a00-4acd-b6db-b7bf1cb1f21b.2b242c02220495576e019c268293a2629917dac4.$locale.xlf", [Type]::Missing, [Type]::Missing, "27f14fa0-3a00-4acd-b6db-b7bf1cb1f21b.2b242c02220495576e019c268293a2629917dac4.$locale.xlf") | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/29f8223291f41775c48fc5f81140839eeff1258d/e2e/91b753c9-5d1c-49cf-8bd3-2fbeee9ef532.md", [Type]::Missing, [Type]::Missing, "91b753c9-5d1c-49cf-8bd3-2fbeee9ef532.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/16eb241e7716f10c06e8c8193ec79918d1a355bd/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/ht/91b753c9-5d1c-49cf-8bd3-2fbeee9ef532.5dd33dce7d1a10de4f7d8751457d7bbec4f4e1ab.$locale.xlf", [Type]::Missing, [Type]::Missing, "91b753c9-5d1c-49cf-8bd3-2fbeee9ef532.5dd33dce7d1a10de4f7d8751457d7bbec4f4e1ab.$locale.xlf") | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/29f8223291f41775c48fc5f81140839eeff1258d/e2e/a6289566-8408-4fec-9a8c-04180367b81e.md", [Type]::Missing, [Type]::Missing, "a6289566-8408-4fec-9a8c-04180367b81e.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/16eb241e7716f10c06e8c8193ec79918d1a355bd/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/ht/a6289566-8408-4fec-9a8c-04180367b81e.12bacf9ab37516007f665f582b00427400306d74.$locale.xlf", [Type]::Missing, [Type]::Missing, "a6289566-8408-4fec-9a8c-04180367b81e.12bacf9ab37516007f665f582b00427400306d74.$locale.xlf") | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/1b714050bd19de299f2bc9b207fc0bdde29e62a9/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null
}

$wsOverview = $wb.Worksheets.Item("Overview")
Rebuild-Overview-Sheet $wsOverview

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Rebuild-Locale-Sheet $wsZhCn "zh-cn" "2016-03-09 20:33:08" "2016-03-09 20:33:41" "2016-03-09 20:34:26" "2016-03-09 20:34:26" "2016-03-09 20:32:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
Rebuild-Locale-Sheet $wsDeDe "de-de" "2016-03-09 20:33:13" "2016-03-09 20:33:53" "2016-03-09 20:34:31" "2016-03-09 20:34:31" "2016-03-09 20:32:41"
